$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1784386617100372
$ws.Range("C2").Value = 0.6059479553903345
$ws.Range("J2").Value = 0.007434944237918215
$ws.Range("P2").Value = 0.1301115241635688
$ws.Range("S2").Value = 0.07806691449814127
$ws.Range("B3").Value = 0.005952380952380952
$ws.Range("C3").Value = 0.01785714285714286
$ws.Range("J3").Value = 0.02380952380952381
$ws.Range("P3").Value = 0.6785714285714286
$ws.Range("S3").Value = 0.2738095238095238
$ws.Range("J4").Value = 0.0196078431372549
$ws.Range("P4").Value = 0.6862745098039216
$ws.Range("S4").Value = 0.2941176470588235
$ws.Range("B6").Value = 0.05343511450381679
$ws.Range("D6").Value = 0.02290076335877863
$ws.Range("F6").Value = 0.07633587786259542
$ws.Range("J6").Value = 0.2480916030534351
$ws.Range("O6").Value = 0.003816793893129771
$ws.Range("Q6").Value = 0.1564885496183206
$ws.Range("R6").Value = 0.04961832061068702
$ws.Range("S6").Value = 0.3893129770992366
$ws.Range("B7").Value = 0.1141304347826087
$ws.Range("D7").Value = 0.02173913043478261
$ws.Range("E7").Value = 0.005434782608695652
$ws.Range("F7").Value = 0.05434782608695652
$ws.Range("J7").Value = 0.1141304347826087
$ws.Range("O7").Value = 0.0108695652173913
$ws.Range("Q7").Value = 0.1358695652173913
$ws.Range("R7").Value = 0.1032608695652174
$ws.Range("S7").Value = 0.4402173913043478
$ws.Range("B8").Value = 0.09223300970873786
$ws.Range("D8").Value = 0.01699029126213592
$ws.Range("E8").Value = 0.002427184466019417
$ws.Range("F8").Value = 0.09223300970873786
$ws.Range("J8").Value = 0.09951456310679611
$ws.Range("O8").Value = 0.01213592233009709
$ws.Range("Q8").Value = 0.2111650485436893
$ws.Range("R8").Value = 0.0703883495145631
$ws.Range("S8").Value = 0.4029126213592233
$ws.Range("B9").Value = 0.07749077490774908
$ws.Range("D9").Value = 0.007380073800738007
$ws.Range("E9").Value = 0.003690036900369004
$ws.Range("F9").Value = 0.07011070110701106
$ws.Range("J9").Value = 0.0959409594095941
$ws.Range("O9").Value = 0.01107011070110701
$ws.Range("Q9").Value = 0.1918819188191882
$ws.Range("R9").Value = 0.0996309963099631
$ws.Range("S9").Value = 0.4428044280442804
$ws.Range("B10").Value = 0.1010971786833856
$ws.Range("D10").Value = 0.0274294670846395
$ws.Range("F10").Value = 0.08072100313479624
$ws.Range("J10").Value = 0.09796238244514106
$ws.Range("O10").Value = 0.009404388714733543
$ws.Range("Q10").Value = 0.2225705329153605
$ws.Range("R10").Value = 0.1073667711598746
$ws.Range("S10").Value = 0.353448275862069
$ws.Range("G11").Value = 0.1717557251908397
$ws.Range("J11").Value = 0.05725190839694656
$ws.Range("K11").Value = 0.2022900763358779
$ws.Range("L11").Value = 0.5648854961832062
$ws.Range("S11").Value = 0.003816793893129771
$ws.Range("G12").Value = 0.74375
$ws.Range("J12").Value = 0.175
$ws.Range("K12").Value = 0.0125
$ws.Range("L12").Value = 0.03125
$ws.Range("S12").Value = 0.0375
$ws.Range("G13").Value = 0.6857142857142857
$ws.Range("J13").Value = 0.2571428571428571
$ws.Range("S13").Value = 0.05714285714285714
$ws.Range("F15").Value = 0.03448275862068965
$ws.Range("H15").Value = 0.1182266009852217
$ws.Range("I15").Value = 0.1083743842364532
$ws.Range("J15").Value = 0.4088669950738916
$ws.Range("K15").Value = 0.0541871921182266
$ws.Range("O15").Value = 0.06403940886699508
$ws.Range("S15").Value = 0.2118226600985222
$ws.Range("F16").Value = 0.02285714285714286
$ws.Range("H16").Value = 0.1485714285714286
$ws.Range("I16").Value = 0.08571428571428572
$ws.Range("K16").Value = 0.1657142857142857
$ws.Range("O16").Value = 0.08
$ws.Range("S16").Value = 0.09714285714285714
$ws.Range("F17").Value = 0.008113590263691683
$ws.Range("H17").Value = 0.18052738336714
$ws.Range("I17").Value = 0.1176470588235294
$ws.Range("J17").Value = 0.4300202839756592
$ws.Range("K17").Value = 0.07099391480730223
$ws.Range("M17").Value = 0.01217038539553753
$ws.Range("N17").Value = 0.002028397565922921
$ws.Range("O17").Value = 0.06288032454361055
$ws.Range("S17").Value = 0.1156186612576065
$ws.Range("F18").Value = 0.0128755364806867
$ws.Range("H18").Value = 0.167381974248927
$ws.Range("I18").Value = 0.1158798283261803
$ws.Range("J18").Value = 0.3991416309012876
$ws.Range("K18").Value = 0.06008583690987124
$ws.Range("M18").Value = 0.02575107296137339
$ws.Range("O18").Value = 0.0815450643776824
$ws.Range("S18").Value = 0.1373390557939914
$ws.Range("F19").Value = 0.02286585365853658
$ws.Range("H19").Value = 0.1745426829268293
$ws.Range("I19").Value = 0.1112804878048781
$ws.Range("J19").Value = 0.3810975609756098
$ws.Range("K19").Value = 0.0861280487804878
$ws.Range("M19").Value = 0.01753048780487805
$ws.Range("N19").Value = 0.001524390243902439
$ws.Range("O19").Value = 0.0586890243902439
$ws.Range("S19").Value = 0.1463414634146341
